# Update timestamp in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 9 de Octubre de 2020 a las 20:50"

# Update country stats (data refresh + re-sort by "Casos totales" descending)
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 7863294
$ws.Cells.Item(4, 3).Value = 29531
$ws.Cells.Item(4, 4).Value = 5041868
$ws.Cells.Item(4, 5).Value = 2603286
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 402
$ws.Cells.Item(4, 8).Value = 218140

$ws.Cells.Item(5, 1).Value = "India"
$ws.Cells.Item(5, 2).Value = 6964074
$ws.Cells.Item(5, 3).Value = 60262
$ws.Cells.Item(5, 4).Value = 5966938
$ws.Cells.Item(5, 5).Value = 889994
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 621
$ws.Cells.Item(5, 8).Value = 107142

$ws.Cells.Item(13, 1).Value = "Francia"
$ws.Cells.Item(13, 2).Value = 691977
$ws.Cells.Item(13, 3).Value = 20339
$ws.Cells.Item(13, 4).Value = 100828
$ws.Cells.Item(13, 5).Value = 558566
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 62
$ws.Cells.Item(13, 8).Value = 32583

$ws.Cells.Item(14, 1).Value = "Sudafrica"
$ws.Cells.Item(14, 2).Value = 686891
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = 618771
$ws.Cells.Item(14, 5).Value = 50712
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 17408

$ws.Cells.Item(25, 1).Value = "Alemania"
$ws.Cells.Item(25, 2).Value = 319403
$ws.Cells.Item(25, 3).Value = 3889
$ws.Cells.Item(25, 4).Value = 269500
$ws.Cells.Item(25, 5).Value = 40227
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 9
$ws.Cells.Item(25, 8).Value = 9676

$ws.Cells.Item(29, 1).Value = "Canada"
$ws.Cells.Item(29, 2).Value = 177697
$ws.Cells.Item(29, 3).Value = 2138
$ws.Cells.Item(29, 4).Value = 149244
$ws.Cells.Item(29, 5).Value = 18867
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(29, 7).Value = 29
$ws.Cells.Item(29, 8).Value = 9586

$ws.Cells.Item(32, 1).Value = "Marruecos"
$ws.Cells.Item(32, 2).Value = 146398
$ws.Cells.Item(32, 3).Value = 3445
$ws.Cells.Item(32, 4).Value = 123022
$ws.Cells.Item(32, 5).Value = 20846
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(32, 7).Value = 44
$ws.Cells.Item(32, 8).Value = 2530

$ws.Cells.Item(33, 1).Value = "Ecuador"
$ws.Cells.Item(33, 2).Value = 145045
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 120511
$ws.Cells.Item(33, 5).Value = 12393
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = 12141

$ws.Cells.Item(34, 1).Value = "Belgica"
$ws.Cells.Item(34, 2).Value = 143596
$ws.Cells.Item(34, 3).Value = 5728
$ws.Cells.Item(34, 4).Value = 19981
$ws.Cells.Item(34, 5).Value = 113489
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 18
$ws.Cells.Item(34, 8).Value = 10126

$ws.Cells.Item(66, 1).Value = "Libano"
$ws.Cells.Item(66, 2).Value = 51170
$ws.Cells.Item(66, 3).Value = 1426
$ws.Cells.Item(66, 4).Value = 22407
$ws.Cells.Item(66, 5).Value = 28313
$ws.Cells.Item(66, 6).Value = 0
$ws.Cells.Item(66, 7).Value = 11
$ws.Cells.Item(66, 8).Value = 450

$ws.Cells.Item(105, 1).Value = "Maldivas"
$ws.Cells.Item(105, 2).Value = 10808
$ws.Cells.Item(105, 3).Value = 66
$ws.Cells.Item(105, 4).Value = 9654
$ws.Cells.Item(105, 5).Value = 1120
$ws.Cells.Item(105, 6).Value = 0
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = 34

$ws.Cells.Item(106, 1).Value = "Georgia"
$ws.Cells.Item(106, 2).Value = 10752
$ws.Cells.Item(106, 3).Value = 527
$ws.Cells.Item(106, 4).Value = 5866
$ws.Cells.Item(106, 5).Value = 4814
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 6
$ws.Cells.Item(106, 8).Value = 72

$ws.Cells.Item(114, 1).Value = "Zimbabue"
$ws.Cells.Item(114, 2).Value = 7994
$ws.Cells.Item(114, 3).Value = 43
$ws.Cells.Item(114, 4).Value = 6474
$ws.Cells.Item(114, 5).Value = 1291
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 229

$ws.Cells.Item(124, 1).Value = "Suazilandia"
$ws.Cells.Item(124, 2).Value = 5644
$ws.Cells.Item(124, 3).Value = 12
$ws.Cells.Item(124, 4).Value = 5254
$ws.Cells.Item(124, 5).Value = 277
$ws.Cells.Item(124, 6).Value = 0
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 113

$ws.Cells.Item(136, 1).Value = "Sri Lanka"
$ws.Cells.Item(136, 2).Value = 4523
$ws.Cells.Item(136, 3).Value = 35
$ws.Cells.Item(136, 4).Value = 3296
$ws.Cells.Item(136, 5).Value = 1214
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 13

$ws.Cells.Item(189, 1).Value = "Monaco"
$ws.Cells.Item(189, 2).Value = 233
$ws.Cells.Item(189, 3).Value = 4
$ws.Cells.Item(189, 4).Value = 209
$ws.Cells.Item(189, 5).Value = 22
$ws.Cells.Item(189, 6).Value = 0
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 2
